$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting the existing weekly records
# (rows 12-137) down by one (to rows 13-138) and growing the used range
# from A1:R137 to A1:R138.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's record. The
# descriptive/constant columns (market, region, category, etc.) match
# every other row in the sheet.
$ws.Cells.Item(12, 1).Value = 3
$ws.Cells.Item(12, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44503
$ws.Cells.Item(12, 5).Value = 5
$ws.Cells.Item(12, 6).Value = 100112010
$ws.Cells.Item(12, 7).Value = "Achicoria"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 110
$ws.Cells.Item(12, 11).Value = 5800
$ws.Cells.Item(12, 12).Value = 6000
$ws.Cells.Item(12, 13).Value = 5891
$ws.Cells.Item(12, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(12, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(12, 16).Value = 368
$ws.Cells.Item(12, 17).Value = 16
$ws.Cells.Item(12, 18).Value = "Hortaliza"
